$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition): update "想去人数" (want-to-go count) values
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 385
$wsExhibition.Range("F3").Value = 2159

# Sheet "全部类型" (All Types): same two events appear again, update accordingly
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 385
$wsAll.Range("F7").Value = 2159
